$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 33266
$ws.Range("E2").Value = -769
$ws.Range("F2").Value = -769
$ws.Range("G2").Value = -507
$ws.Range("H2").Value = -509
$ws.Range("I2").Value = -509
$ws.Range("K2").Value = 19684
$ws.Range("L2").Value = 11910
$ws.Range("M2").Value = 7774
$ws.Range("N2").Value = 7774
$ws.Range("P2").Value = 6861
$ws.Range("Q2").Value = 682
$ws.Range("R2").Value = -2575
$ws.Range("S2").Value = -214
$ws.Range("T2").Value = 1837
$ws.Range("U2").Value = -1155
$ws.Range("V2").Value = 1245
$ws.Range("W2").Value = -2.31
$ws.Range("X2").Value = -1.53
$ws.Range("Y2").Value = -6.13
$ws.Range("Z2").Value = -2.49
$ws.Range("AA2").Value = 153.21
$ws.Range("AB2").Value = 15.25
$ws.Range("AC2").Value = -371
$ws.Range("AD2").Value = -24.95
$ws.Range("AE2").Value = 5666
$ws.Range("AF2").Value = 1.63
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 137220096
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 33901
$ws.Range("E3").Value = -358
$ws.Range("F3").Value = -358
$ws.Range("G3").Value = -619
$ws.Range("H3").Value = -619
$ws.Range("I3").Value = -619
$ws.Range("K3").Value = 20392
$ws.Range("L3").Value = 13221
$ws.Range("M3").Value = 7171
$ws.Range("N3").Value = 7171
$ws.Range("P3").Value = 6861
$ws.Range("Q3").Value = 2015
$ws.Range("R3").Value = -2280
$ws.Range("S3").Value = 675
$ws.Range("T3").Value = 1511
$ws.Range("U3").Value = 505
$ws.Range("V3").Value = 1907
$ws.Range("W3").Value = -1.06
$ws.Range("X3").Value = -1.83
$ws.Range("Y3").Value = -8.289999999999999
$ws.Range("Z3").Value = -3.09
$ws.Range("AA3").Value = 184.38
$ws.Range("AB3").Value = 4.15
$ws.Range("AC3").Value = -451
$ws.Range("AD3").Value = -17.17
$ws.Range("AE3").Value = 5226
$ws.Range("AF3").Value = 1.48
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 137220096
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 36285
$ws.Range("E4").Value = 280
$ws.Range("F4").Value = 280
$ws.Range("G4").Value = 581
$ws.Range("H4").Value = 581
$ws.Range("I4").Value = 581
$ws.Range("K4").Value = 21634
$ws.Range("L4").Value = 13507
$ws.Range("M4").Value = 8127
$ws.Range("N4").Value = 8127
$ws.Range("P4").Value = 6861
$ws.Range("Q4").Value = 2444
$ws.Range("R4").Value = -2107
$ws.Range("S4").Value = 70
$ws.Range("T4").Value = 1275
$ws.Range("U4").Value = 1169
$ws.Range("V4").Value = 1945
$ws.Range("W4").Value = 0.77
$ws.Range("X4").Value = 1.6
$ws.Range("Y4").Value = 7.6
$ws.Range("Z4").Value = 2.77
$ws.Range("AA4").Value = 166.2
$ws.Range("AB4").Value = 18.5
$ws.Range("AC4").Value = 423
$ws.Range("AD4").Value = 18.63
$ws.Range("AE4").Value = 5923
$ws.Range("AF4").Value = 1.33
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 137220096
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 34946
$ws.Range("E5").Value = -653
$ws.Range("F5").Value = -653
$ws.Range("G5").Value = -658
$ws.Range("H5").Value = -658
$ws.Range("I5").Value = -658
$ws.Range("K5").Value = 22480
$ws.Range("L5").Value = 14729
$ws.Range("M5").Value = 7752
$ws.Range("N5").Value = 7752
$ws.Range("P5").Value = 6897
$ws.Range("Q5").Value = 2044
$ws.Range("R5").Value = -2775
$ws.Range("S5").Value = 502
$ws.Range("T5").Value = 1575
$ws.Range("U5").Value = 469
$ws.Range("V5").Value = 2338
$ws.Range("W5").Value = -1.87
$ws.Range("X5").Value = -1.88
$ws.Range("Y5").Value = -8.289999999999999
$ws.Range("Z5").Value = -2.98
$ws.Range("AA5").Value = 190.01
$ws.Range("AB5").Value = 12.08
$ws.Range("AC5").Value = -478
$ws.Range("AD5").Value = -10.72
$ws.Range("AE5").Value = 5619
$ws.Range("AF5").Value = 0.91
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 137949396
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 37048
$ws.Range("E6").Value = -642
$ws.Range("F6").Value = -642
$ws.Range("G6").Value = -618
$ws.Range("H6").Value = -618
$ws.Range("I6").Value = -618
$ws.Range("K6").Value = 22056
$ws.Range("L6").Value = 15124
$ws.Range("M6").Value = 6933
$ws.Range("N6").Value = 6933
$ws.Range("P6").Value = 6897
$ws.Range("Q6").Value = 1771
$ws.Range("R6").Value = -2577
$ws.Range("S6").Value = 213
$ws.Range("T6").Value = 1369
$ws.Range("U6").Value = 402
$ws.Range("V6").Value = 2532
$ws.Range("W6").Value = -1.73
$ws.Range("X6").Value = -1.67
$ws.Range("Y6").Value = -8.42
$ws.Range("Z6").Value = -2.78
$ws.Range("AA6").Value = 218.14
$ws.Range("AB6").Value = 0.2
$ws.Range("AC6").Value = -448
$ws.Range("AD6").Value = -8.83
$ws.Range("AE6").Value = 5026
$ws.Range("AF6").Value = 0.79
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 137949396
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 36680
$ws.Range("E7").Value = -1830
$ws.Range("G7").Value = -1810
$ws.Range("H7").Value = -1810
$ws.Range("I7").Value = -1810
$ws.Range("K7").Value = 21080
$ws.Range("L7").Value = 15960
$ws.Range("M7").Value = 5120
$ws.Range("N7").Value = 5120
$ws.Range("P7").Value = 6900
$ws.Range("Q7").Value = 1420
$ws.Range("R7").Value = -2580
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 1390
$ws.Range("U7").Value = -1160
$ws.Range("W7").Value = -4.99
$ws.Range("X7").Value = -4.93
$ws.Range("Y7").Value = -30.03
$ws.Range("Z7").Value = -8.390000000000001
$ws.Range("AA7").Value = 311.72
$ws.Range("AC7").Value = -1214
$ws.Range("AD7").Value = -1.61
$ws.Range("AE7").Value = 3417
$ws.Range("AF7").Value = 0.57
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 38140
$ws.Range("E8").Value = 0
$ws.Range("G8").Value = -10
$ws.Range("H8").Value = -10
$ws.Range("I8").Value = -10
$ws.Range("K8").Value = 21700
$ws.Range("L8").Value = 16580
$ws.Range("M8").Value = 5120
$ws.Range("N8").Value = 5120
$ws.Range("P8").Value = 6900
$ws.Range("Q8").Value = 3450
$ws.Range("R8").Value = -2620
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 1410
$ws.Range("U8").Value = 830
$ws.Range("X8").Value = -0.03
$ws.Range("Y8").Value = -0.19
$ws.Range("Z8").Value = -0.05
$ws.Range("AA8").Value = 323.83
$ws.Range("AC8").Value = -7
$ws.Range("AD8").Value = -292.19
$ws.Range("AE8").Value = 3417
$ws.Range("AF8").Value = 0.57
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("W8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 38910
$ws.Range("E9").Value = 390
$ws.Range("G9").Value = 360
$ws.Range("H9").Value = 360
$ws.Range("I9").Value = 360
$ws.Range("K9").Value = 22250
$ws.Range("L9").Value = 16770
$ws.Range("M9").Value = 5480
$ws.Range("N9").Value = 5480
$ws.Range("P9").Value = 6900
$ws.Range("Q9").Value = 3030
$ws.Range("R9").Value = -2620
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 1400
$ws.Range("U9").Value = 420
$ws.Range("W9").Value = 1
$ws.Range("X9").Value = 0.93
$ws.Range("Y9").Value = 6.79
$ws.Range("Z9").Value = 1.64
$ws.Range("AA9").Value = 306.02
$ws.Range("AC9").Value = 240
$ws.Range("AD9").Value = 8.119999999999999
$ws.Range("AE9").Value = 3657
$ws.Range("AF9").Value = 0.53
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
